# Auto-generated edit script: apply cell value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 42784664
$ws.Range("I116").Value = 42784664
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 42784664
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -42781222
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 26317106
$ws.Range("I137").Value = 34483684
$ws.Range("J137").Value = 2571.7778
$ws.Range("K137").Value = 103451052
$ws.Range("L137").Value = 7715.3334
$ws.Range("M137").Value = -103448502
$ws.Range("N137").Value = -12815.3334
$ws.Range("H138").Value = 7513944
$ws.Range("I138").Value = 1450148.6
$ws.Range("J138").Value = 10419513
$ws.Range("K138").Value = 4350445.800000001
$ws.Range("L138").Value = 31258539
$ws.Range("M138").Value = -4345305.800000001
$ws.Range("N138").Value = -31268819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2524.92
$ws.Range("I61").Value = 1827.6666
$ws.Range("J61").Value = 3878.4119
$ws.Range("K61").Value = 1827.6666
$ws.Range("L61").Value = 3878.4119
$ws.Range("M61").Value = -1615.6666
$ws.Range("N61").Value = -4302.4119
$ws.Range("H74").Value = 5383.0967
$ws.Range("I74").Value = 1504.174
$ws.Range("J74").Value = 16535
$ws.Range("K74").Value = 1504.174
$ws.Range("L74").Value = 16535
$ws.Range("M74").Value = -630.174
$ws.Range("N74").Value = -18283
$ws.Range("H77").Value = 5383.0967
$ws.Range("I77").Value = 1504.174
$ws.Range("J77").Value = 16535
$ws.Range("K77").Value = 7520.87
$ws.Range("L77").Value = 82675
$ws.Range("M77").Value = -3152.87
$ws.Range("N77").Value = -91411
$ws.Range("H122").Value = 990.6875
$ws.Range("I122").Value = 877
$ws.Range("J122").Value = 1331.75
$ws.Range("K122").Value = 2631
$ws.Range("L122").Value = 3995.25
$ws.Range("M122").Value = -181
$ws.Range("N122").Value = -8895.25
$ws.Range("H132").Value = 2321.3618
$ws.Range("I132").Value = 1838.5758
$ws.Range("J132").Value = 3459.3572
$ws.Range("K132").Value = 5515.7274
$ws.Range("L132").Value = 10378.0716
$ws.Range("M132").Value = -2985.7274
$ws.Range("N132").Value = -15438.0716
$ws.Range("H136").Value = 2524.92
$ws.Range("I136").Value = 1827.6666
$ws.Range("J136").Value = 3878.4119
$ws.Range("K136").Value = 5482.9998
$ws.Range("L136").Value = 11635.2357
$ws.Range("M136").Value = -2932.9998
$ws.Range("N136").Value = -16735.2357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4390.857
$ws.Range("I86").Value = 1283.5333
$ws.Range("J86").Value = 6721.35
$ws.Range("K86").Value = 1283.5333
$ws.Range("L86").Value = 6721.35
$ws.Range("M86").Value = -160.5333000000001
$ws.Range("N86").Value = -8967.35
$ws.Range("H89").Value = 4390.857
$ws.Range("I89").Value = 1283.5333
$ws.Range("J89").Value = 6721.35
$ws.Range("K89").Value = 6417.6665
$ws.Range("L89").Value = 33606.75
$ws.Range("M89").Value = -801.6665000000003
$ws.Range("N89").Value = -44838.75
$ws.Range("H128").Value = 1000
$ws.Range("I128").Value = 1000
$ws.Range("K128").Value = 3000
$ws.Range("M128").Value = -510
$ws.Range("H134").Value = 4847
$ws.Range("I134").Value = 3591
$ws.Range("J134").Value = 5760.4546
$ws.Range("K134").Value = 10773
$ws.Range("L134").Value = 17281.3638
$ws.Range("M134").Value = -8238
$ws.Range("N134").Value = -22351.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1262.1154
$ws.Range("I31").Value = 1192.2916
$ws.Range("J31").Value = 2100
$ws.Range("K31").Value = 1192.2916
$ws.Range("L31").Value = 2100
$ws.Range("M31").Value = -897.2916
$ws.Range("N31").Value = -2690
$ws.Range("H34").Value = 1262.1154
$ws.Range("I34").Value = 1192.2916
$ws.Range("J34").Value = 2100
$ws.Range("K34").Value = 1192.2916
$ws.Range("L34").Value = 2100
$ws.Range("M34").Value = -990.2916
$ws.Range("N34").Value = -2504
$ws.Range("H58").Value = 1786.9412
$ws.Range("I58").Value = 1334.2174
$ws.Range("J58").Value = 2733.5454
$ws.Range("K58").Value = 1334.2174
$ws.Range("L58").Value = 2733.5454
$ws.Range("M58").Value = -1131.2174
$ws.Range("N58").Value = -3139.5454
$ws.Range("H132").Value = 2275.535
$ws.Range("I132").Value = 1661.9032
$ws.Range("J132").Value = 3860.75
$ws.Range("K132").Value = 4985.7096
$ws.Range("L132").Value = 11582.25
$ws.Range("M132").Value = -2455.7096
$ws.Range("N132").Value = -16642.25
$ws.Range("H134").Value = 2602.0715
$ws.Range("I134").Value = 691.63635
$ws.Range("J134").Value = 9607
$ws.Range("K134").Value = 2074.90905
$ws.Range("L134").Value = 28821
$ws.Range("M134").Value = 460.0909499999998
$ws.Range("N134").Value = -33891
$ws.Range("H136").Value = 1786.9412
$ws.Range("I136").Value = 1334.2174
$ws.Range("J136").Value = 2733.5454
$ws.Range("K136").Value = 4002.6522
$ws.Range("L136").Value = 8200.636200000001
$ws.Range("M136").Value = -1452.6522
$ws.Range("N136").Value = -13300.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 412.5
$ws.Range("I98").Value = 383.33334
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 1150.00002
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 347.9999800000001
$ws.Range("N98").Value = -4496
$ws.Range("H105").Value = 7971.6
$ws.Range("J105").Value = 7971.6
$ws.Range("L105").Value = 23914.8
$ws.Range("N105").Value = -29156.8
$ws.Range("H128").Value = 191592
$ws.Range("I128").Value = 191592
$ws.Range("K128").Value = 574776
$ws.Range("M128").Value = -569796

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 108285
$ws.Range("J109").Value = 108285
$ws.Range("L109").Value = 108285
$ws.Range("N109").Value = -110365
$ws.Range("H121").Value = 25982
$ws.Range("J121").Value = 25982
$ws.Range("L121").Value = 25982
$ws.Range("N121").Value = -29476
$ws.Range("H132").Value = 2672.568
$ws.Range("I132").Value = 2216.606
$ws.Range("K132").Value = 6649.818000000001
$ws.Range("M132").Value = -4119.818000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 1040
$ws.Range("I107").Value = 1040
$ws.Range("K107").Value = 1040
$ws.Range("M107").Value = 880
$ws.Range("H132").Value = 4902.3794
$ws.Range("I132").Value = 4041.0908
$ws.Range("J132").Value = 5428.722
$ws.Range("K132").Value = 12123.2724
$ws.Range("L132").Value = 16286.166
$ws.Range("M132").Value = -9593.2724
$ws.Range("N132").Value = -21346.166
$ws.Range("H136").Value = 4832
$ws.Range("I136").Value = 2970.0454
$ws.Range("J136").Value = 8245.583000000001
$ws.Range("K136").Value = 8910.136200000001
$ws.Range("L136").Value = 24736.749
$ws.Range("M136").Value = -6360.136200000001
$ws.Range("N136").Value = -29836.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 420.55554
$ws.Range("I107").Value = 396.72726
$ws.Range("J107").Value = 458
$ws.Range("K107").Value = 1190.18178
$ws.Range("L107").Value = 1374
$ws.Range("M107").Value = 729.8182200000001
$ws.Range("N107").Value = -5214
$ws.Range("H132").Value = 9806171
$ws.Range("I132").Value = 11906542
$ws.Range("K132").Value = 35719626
$ws.Range("M132").Value = -35717096
$ws.Range("H136").Value = 12384780
$ws.Range("I136").Value = 23882118
$ws.Range("J136").Value = 3031.2307
$ws.Range("K136").Value = 71646354
$ws.Range("L136").Value = 9093.6921
$ws.Range("M136").Value = -71643804
$ws.Range("N136").Value = -14193.6921
